$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value = 7593
$ws.Range("F12").Value = 8226
$ws.Range("F16").Value = 5629
$ws.Range("F17").Value = 5629
$ws.Range("F19").Value = 2580
$ws.Range("F20").Value = 1116
$ws.Range("F21").Value = 4587
$ws.Range("F22").Value = 335
$ws.Range("F24").Value = 91
$ws.Range("F26").Value = 506
$ws.Range("F27").Value = 3298
$ws.Range("F28").Value = 3298
$ws.Range("F30").Value = 12
$ws.Range("F31").Value = 2876
$ws.Range("F32").Value = 2876
$ws.Range("F35").Value = 122
$ws.Range("F36").Value = 287
$ws.Range("F37").Value = 4
$ws.Range("G40").Value = 9.9
$ws.Range("F41").Value = 1644
$ws.Range("F44").Value = 7
$ws.Range("F45").Value = 2661
$ws.Range("F46").Value = 2
$ws.Range("F47").Value = 2274
$ws.Range("F48").Value = 8
$ws.Range("F49").Value = 27

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 98
$ws.Range("F3").Value = 110

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 7593
$ws.Range("F11").Value = 8226
$ws.Range("F14").Value = 5629
$ws.Range("F15").Value = 5629
$ws.Range("F17").Value = 2580
$ws.Range("F18").Value = 1116
$ws.Range("F19").Value = 4587
$ws.Range("F21").Value = 91
$ws.Range("F22").Value = 98
$ws.Range("F24").Value = 110
$ws.Range("F25").Value = 506
$ws.Range("F26").Value = 3298
$ws.Range("F27").Value = 3298
$ws.Range("F29").Value = 12
$ws.Range("F30").Value = 2876
$ws.Range("F31").Value = 2876
$ws.Range("F33").Value = 122
$ws.Range("F34").Value = 287
$ws.Range("F36").Value = 4
$ws.Range("G40").Value = 9.9
$ws.Range("F42").Value = 1644
$ws.Range("F45").Value = 7
$ws.Range("F46").Value = 2661
$ws.Range("F48").Value = 2274
$ws.Range("F49").Value = 8
$ws.Range("F50").Value = 27
